# Updates cryptos list values (coin name/link reorderings and price/volume refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.385.57"
$ws.Range("E2").Value = "'  +0.69%  "

$ws.Range("D3").Value = "'3.125.91"
$ws.Range("E3").Value = "'  +0.31%  "

$ws.Range("E4").Value = "'  +0.04%  "

$ws.Range("D5").Value = "'241.07"
$ws.Range("E5").Value = "'  -1.28%  "

$ws.Range("D6").Value = "'618.45"
$ws.Range("E6").Value = "'  -1.14%  "

$ws.Range("D7").Value = "'1.11"
$ws.Range("E7").Value = "'  -3.56%  "

$ws.Range("D8").Value = "'0.386"
$ws.Range("E8").Value = "'  +3.63%  "

$ws.Range("E9").Value = "'  -0.05%  "

$ws.Range("D10").Value = "'3.124.86"
$ws.Range("E10").Value = "'  +0.36%  "

$ws.Range("D11").Value = "'0.747"
$ws.Range("E11").Value = "'  -0.25%  "

$ws.Range("D12").Value = "'0.203"

$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "'  -0.43%  "

$ws.Range("D14").Value = "'35.02"
$ws.Range("E14").Value = "'  -1.11%  "

$ws.Range("D15").Value = "'5.60"
$ws.Range("E15").Value = "'  +1.94%  "

$ws.Range("D16").Value = "'91.285.73"
$ws.Range("E16").Value = "'  +0.84%  "

$ws.Range("E17").Value = "'  +0.84%  "

$ws.Range("D18").Value = "'3.169.95"
$ws.Range("E18").Value = "'  +2.31%  "

$ws.Range("D19").Value = "'3.75"
$ws.Range("E19").Value = "'  -1.89%  "

$ws.Range("D20").Value = "'14.92"
$ws.Range("E20").Value = "'  +3.61%  "

$ws.Range("E21").Value = "'  +2.51%  "

$ws.Range("D22").Value = "'455.61"
$ws.Range("E22").Value = "'  +1.66%  "

$ws.Range("D23").Value = "'0.0000201"
$ws.Range("E23").Value = "'  -4.68%  "

$ws.Range("D24").Value = "'9.17"
$ws.Range("E24").Value = "'  +0.37%  "

$ws.Range("D25").Value = "'5.90"
$ws.Range("E25").Value = "'  -0.22%  "

$ws.Range("D26").Value = "'88.69"
$ws.Range("E26").Value = "'  -5.10%  "

$ws.Range("D27").Value = "'11.80"
$ws.Range("E27").Value = "'  -1.53%  "

$ws.Range("D28").Value = "'0.151"
$ws.Range("E28").Value = "'  +36.53%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  -0.12%  "

$ws.Range("D31").Value = "'0.230"
$ws.Range("E31").Value = "'  +4.89%  "

$ws.Range("D32").Value = "'0.166"
$ws.Range("E32").Value = "'  -6.74%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.176"
$ws.Range("E33").Value = "'  +10.67%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'9.33"
$ws.Range("E34").Value = "'  +1.22%  "

$ws.Range("D35").Value = "'26.36"
$ws.Range("E35").Value = "'  -1.03%  "

$ws.Range("D36").Value = "'7.43"
$ws.Range("E36").Value = "'  -2.93%  "

$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = "'  +1.79%  "

$ws.Range("B38").Value = "MantraDAO"
$ws.Range("C38").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D38").Value = "'3.93"
$ws.Range("E38").Value = "'  -6.94%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'490.90"
$ws.Range("E39").Value = "'  -0.80%  "

$ws.Range("E40").Value = "'  +1.52%  "

$ws.Range("D41").Value = "'0.438"
$ws.Range("E41").Value = "'  +4.96%  "

$ws.Range("E42").Value = "'  -5.95%  "

$ws.Range("E43").Value = "'  +0.07%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "'  -0.02%  "

$ws.Range("B45").Value = "Binance-PegBSC-USD"
$ws.Range("C45").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D45").Value = "'0.720"
$ws.Range("E45").Value = "'  -27.88%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'0.706"
$ws.Range("E46").Value = "'  +2.64%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'156.18"
$ws.Range("E47").Value = "'  -0.77%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.92"
$ws.Range("E48").Value = "'  +0.51%  "

$ws.Range("D49").Value = "'1.35"
$ws.Range("E49").Value = "'  +0.64%  "

$ws.Range("E50").Value = "'  -2.11%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "'44.04"
$ws.Range("E51").Value = "'  -2.19%  "
